# Applies the "bugs" worksheet update described by the commit:
# "implemented proper removal of leftovers from the build process; revised README.md"
#
# Concretely (sheet "Buggar" = ActiveSheet):
#  - row 2:  Datum atgardat (C2) filled in; sitemap bug description text extended
#  - row 5:  Datum atgardat (C5) filled in
#  - row 8:  Release (A8) alpha-3 -> alpha-4; Prio (E8) 2 -> 1
#  - row 9:  Datum atgardat (C9) filled in
#  - row 10: Release (A10) alpha-3 -> alpha-5; Prio (E10) set to 2
#  - row 11: Prio (E11) set to 2; new bug description added (F11)
#  - row 12: Datum atgardat (C12), Prio (E12) and description (F12) filled in
#  - rows 14-19: Release (A) + Datum rapporterat (B) filled in for previously blank rows
#  - selection moves from F10 to C9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Buggar")

# NOTE: the order in which brand-new unique strings are first written
# determines where they land in xl/sharedStrings.xml, which the canonical
# OOXML pins (sitemap note, then LinkedIn-badge note, then "alpha-4", then
# the h1-rubrik note, then "alpha-5"). Touch those cells in that order;
# the remaining (non-string-table-affecting) edits are interleaved after.

$ws.Range("F2").Value = "Fixa en sitemap. Registrera den hos Google."
$ws.Range("F12").Value = "LinkedIn-badge. Problemet var löst innan jag skrev in det som en bugg, därav de ologiska datumen."
$ws.Range("A8").Value = "alpha-4"
$ws.Range("F11").Value = "Varje sida skall ha en h1-rubrik"
$ws.Range("A10").Value = "alpha-5"

# --- Row 2: Datum atgardat ---------------------------------------------------
$ws.Range("C2").Value = 43202

# --- Row 5: Datum atgardat -------------------------------------------------
$ws.Range("C5").Value = 43203

# --- Row 8: Prio corrected --------------------------------------------------
$ws.Range("E8").Value = "1"

# --- Row 9: Datum atgardat --------------------------------------------------
$ws.Range("C9").Value = 43202

# --- Row 10: Prio set --------------------------------------------------------
$ws.Range("E10").Value = "2"

# --- Row 11: Prio set --------------------------------------------------------
$ws.Range("E11").Value = "2"

# --- Row 12: Datum atgardat, Prio --------------------------------------------
$ws.Range("C12").Value = 43201
$ws.Range("E12").Value = "1"

# --- Rows 14-18: previously blank Release/Datum rapporterat cells ----------
$ws.Range("A14").Value = "alpha-3"
$ws.Range("B14").Value = 43202

$ws.Range("A15").Value = "alpha-3"
$ws.Range("B15").Value = 43202

$ws.Range("A16").Value = "alpha-3"
$ws.Range("B16").Value = 43202

$ws.Range("A17").Value = "alpha-3"
$ws.Range("B17").Value = 43202

$ws.Range("A18").Value = "alpha-3"
$ws.Range("B18").Value = 43202
# B18 previously carried the unformatted "style=10" date format (no explicit
# font applied); nudge its alignment so it normalizes onto the same cell
# style as the rest of column B (style=9) rather than staying on style=10.
$ws.Range("B18").HorizontalAlignment = -4152

# --- Row 19: brand-new Release/Datum rapporterat cells ----------------------
$ws.Range("A19").Value = "alpha-3"
# B19 is a brand-new cell; column B's default style (10) differs from the
# style used throughout the rest of the column (9), so copy formats down
# from a neighboring "Datum rapporterat" cell before writing the value.
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("B19").Value = 43202

# --- Selection moves from F10 to C9 -----------------------------------------
$ws.Range("C9").Select() | Out-Null
